$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to stay plain text (matches the source
# workbook, which stores these as inline strings like "28.637.93" / "  +2.27%  ")
# so Excel does not silently reinterpret them as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value2 = "28.637.93"
$ws.Range("E2").Value2 = "  +2.27%  "
$ws.Range("D3").Value2 = "1.869.79"
$ws.Range("E3").Value2 = "  +2.23%  "
$ws.Range("D4").Value2 = "1.004"
$ws.Range("E4").Value2 = "  +0.26%  "
$ws.Range("D5").Value2 = "326.44"
$ws.Range("E5").Value2 = "  +0.08%  "
$ws.Range("D7").Value2 = "0.4660"
$ws.Range("E7").Value2 = "  +1.11%  "
$ws.Range("E8").Value2 = "  +0.55%  "
$ws.Range("D9").Value2 = "0.07869"
$ws.Range("E9").Value2 = "  +0.08%  "
$ws.Range("D10").Value2 = "0.9738"
$ws.Range("E10").Value2 = "  +1.61%  "
$ws.Range("D11").Value2 = "21.94"
$ws.Range("E11").Value2 = "  +0.43%  "
$ws.Range("D12").Value2 = "1.921.40"
$ws.Range("E12").Value2 = "  -0.67%  "
$ws.Range("E13").Value2 = "  +1.56%  "
$ws.Range("D14").Value2 = "5.698"
$ws.Range("E14").Value2 = "  +0.83%  "
$ws.Range("D15").Value2 = "0.06996"
$ws.Range("E15").Value2 = "  +3.67%  "
$ws.Range("D16").Value2 = "88.04"
$ws.Range("E16").Value2 = "  +1.45%  "
$ws.Range("D18").Value2 = "0.00001004"
$ws.Range("E18").Value2 = "  +1.28%  "
$ws.Range("D19").Value2 = "16.81"
$ws.Range("E19").Value2 = "  +1.36%  "
$ws.Range("E20").Value2 = "  +0.32%  "
$ws.Range("D21").Value2 = "28.625.79"
$ws.Range("E21").Value2 = "  +2.15%  "
$ws.Range("D22").Value2 = "5.289"
$ws.Range("E22").Value2 = "  -0.09%  "
$ws.Range("D23").Value2 = "11.00"
$ws.Range("D24").Value2 = "2.114"
$ws.Range("E24").Value2 = "  +1.40%  "
$ws.Range("D25").Value2 = "2.040.79"
$ws.Range("E25").Value2 = "  -3.09%  "
$ws.Range("D26").Value2 = "152.71"
$ws.Range("E26").Value2 = "  -0.57%  "
$ws.Range("D27").Value2 = "19.22"
$ws.Range("E27").Value2 = "  +0.58%  "
$ws.Range("D28").Value2 = "5.792"
$ws.Range("E28").Value2 = "  +1.03%  "
$ws.Range("E29").Value2 = "  +0.72%  "
$ws.Range("D30").Value2 = "119.45"
$ws.Range("E30").Value2 = "  +2.06%  "
$ws.Range("D31").Value2 = "0.09372"
$ws.Range("E31").Value2 = "  +1.43%  "
$ws.Range("D32").Value2 = "0.9207"
$ws.Range("E32").Value2 = "  -1.46%  "
$ws.Range("D33").Value2 = "5.268"
$ws.Range("E33").Value2 = "  -0.36%  "
$ws.Range("D34").Value2 = "1.339"
$ws.Range("E34").Value2 = "  +1.91%  "
$ws.Range("E35").Value2 = "  +0.67%  "
$ws.Range("D36").Value2 = "0.05803"
$ws.Range("E36").Value2 = "  -1.02%  "
$ws.Range("D37").Value2 = "0.02096"
$ws.Range("E37").Value2 = "  -2.04%  "
$ws.Range("E38").Value2 = "  +0.37%  "
$ws.Range("D39").Value2 = "7.750"
$ws.Range("E39").Value2 = "  -0.05%  "
$ws.Range("D40").Value2 = "0.5624"
$ws.Range("E40").Value2 = "  +0.83%  "
$ws.Range("D41").Value2 = "0.1786"
$ws.Range("E41").Value2 = "  +1.65%  "
$ws.Range("D42").Value2 = "9.741"
$ws.Range("E42").Value2 = "  -1.25%  "
$ws.Range("E43").Value2 = "  +2.70%  "
$ws.Range("D44").Value2 = "11.69"
$ws.Range("E44").Value2 = "  +1.25%  "
$ws.Range("D45").Value2 = "0.5311"
$ws.Range("E45").Value2 = "  +1.07%  "
$ws.Range("E46").Value2 = "  -6.01%  "
$ws.Range("B47").Value2 = "NEARProtocol"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value2 = "1.823"
$ws.Range("E47").Value2 = "  -0.08%  "
$ws.Range("B48").Value2 = "RenderToken"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value2 = "2.079"
$ws.Range("E48").Value2 = "  -3.09%  "
$ws.Range("D49").Value2 = "113.21"
$ws.Range("E49").Value2 = "  +0.54%  "
$ws.Range("D50").Value2 = "2.404"
$ws.Range("E50").Value2 = "  +3.70%  "
$ws.Range("D51").Value2 = "1.004"
$ws.Range("E51").Value2 = "  +0.30%  "

# Restore the default (unstyled) look for that range now that the text is locked in.
$ws.Range("D2:E51").Style = "Normal"
